$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 "Save", matching the style of the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column H for rows 2-4
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 1
